# Updated RAD Test Cases and data to handle FEIN/SSN and Failures when
# Payment Applications are not deployed in QA2.
#
# Each worksheet is a Katalon RAD test-case log: column A = Result
# (Pass/Fail), column B = execution Date/time, columns C.. = the test's
# input parameters. This script stamps a fresh run's Result/Date into
# column A/B for every test row across the eight sheets.
#
# New, never-before-seen cells are nudged back to the "no explicit style"
# state (matching their column-A/B siblings already on the sheet) by
# resetting them to the Normal cell style right after the value is
# written - otherwise a brand new cell inherits the column's style index.

$wb = $excel.ActiveWorkbook

function Set-Result {
    param($ws, [int]$row, [string]$result, [string]$date)

    $cellA = $ws.Cells.Item($row, 1)
    $wasEmptyA = ($cellA.Value -eq $null)
    $cellA.Value = $result
    if (-not $wasEmptyA) {
        $cellA.Style = "Normal"
    }

    $cellB = $ws.Cells.Item($row, 2)
    $wasEmptyB = ($cellB.Value -eq $null)
    $cellB.Value = $date
    if (-not $wasEmptyB) {
        $cellB.Style = "Normal"
    }
}

# ---- Estimated (sheet1) ---------------------------------------------
# Row 3 is written before row 2 so the new shared strings land in the
# same order as the reference workbook ("Fail" + its date before the
# second "Pass" date).
$ws1 = $wb.Worksheets.Item("Estimated")
Set-Result $ws1 3 "Fail" "Fri Sep 29 16:26:03 EDT 2023"
Set-Result $ws1 2 "Pass" "Fri Sep 29 16:27:50 EDT 2023"

# ---- Existing (sheet2) -----------------------------------------------
$ws2 = $wb.Worksheets.Item("Existing")
Set-Result $ws2 2  "Pass" "Fri Sep 29 16:38:16 EDT 2023"
Set-Result $ws2 3  "Pass" "Fri Sep 29 16:39:04 EDT 2023"
Set-Result $ws2 4  "Pass" "Fri Sep 29 16:39:48 EDT 2023"
Set-Result $ws2 5  "Pass" "Fri Sep 29 16:40:34 EDT 2023"
Set-Result $ws2 6  "Pass" "Fri Sep 29 16:41:20 EDT 2023"
Set-Result $ws2 7  "Pass" "Fri Sep 29 16:42:04 EDT 2023"
Set-Result $ws2 8  "Pass" "Fri Sep 29 16:42:49 EDT 2023"
Set-Result $ws2 9  "Pass" "Fri Sep 29 16:43:33 EDT 2023"
Set-Result $ws2 10 "Pass" "Fri Sep 29 16:44:18 EDT 2023"
Set-Result $ws2 11 "Fail" "Fri Sep 29 16:45:01 EDT 2023"
Set-Result $ws2 12 "Fail" "Fri Sep 29 16:46:01 EDT 2023"

# ---- Extension (sheet3) ----------------------------------------------
$ws3 = $wb.Worksheets.Item("Extension")
Set-Result $ws3 2 "Pass" "Fri Sep 29 16:56:41 EDT 2023"
Set-Result $ws3 3 "Fail" "Fri Sep 29 16:57:27 EDT 2023"
Set-Result $ws3 4 "Fail" "Fri Sep 29 16:58:26 EDT 2023"
Set-Result $ws3 5 "Fail" "Fri Sep 29 16:59:25 EDT 2023"
Set-Result $ws3 6 "Fail" "Fri Sep 29 17:00:24 EDT 2023"
Set-Result $ws3 7 "Fail" "Fri Sep 29 17:01:24 EDT 2023"

# ---- NewTaxReturn (sheet4) --------------------------------------------
$ws4 = $wb.Worksheets.Item("NewTaxReturn")
Set-Result $ws4 2  "Fail" "Fri Sep 29 17:08:47 EDT 2023"
Set-Result $ws4 3  "Fail" "Fri Sep 29 17:09:49 EDT 2023"
Set-Result $ws4 4  "Fail" "Fri Sep 29 17:10:47 EDT 2023"
Set-Result $ws4 5  "Fail" "Fri Sep 29 17:11:46 EDT 2023"
Set-Result $ws4 6  "Fail" "Fri Sep 29 17:12:45 EDT 2023"
Set-Result $ws4 7  "Fail" "Fri Sep 29 17:13:43 EDT 2023"
Set-Result $ws4 8  "Fail" "Fri Sep 29 17:14:42 EDT 2023"
Set-Result $ws4 9  "Fail" "Fri Sep 29 17:15:40 EDT 2023"
Set-Result $ws4 10 "Fail" "Fri Sep 29 17:16:38 EDT 2023"
Set-Result $ws4 11 "Fail" "Fri Sep 29 17:17:37 EDT 2023"
Set-Result $ws4 12 "Fail" "Fri Sep 29 17:18:36 EDT 2023"
Set-Result $ws4 13 "Fail" "Fri Sep 29 17:19:34 EDT 2023"
Set-Result $ws4 14 "Fail" "Fri Sep 29 17:20:33 EDT 2023"
Set-Result $ws4 15 "Fail" "Fri Sep 29 17:21:32 EDT 2023"
Set-Result $ws4 16 "Fail" "Fri Sep 29 17:22:30 EDT 2023"

# ---- Personal_EL (sheet8) --------------------------------------------
# Written before Personal_IND/Personal_JNT so the new shared string for
# this run's date lands immediately after NewTaxReturn's, matching the
# reference workbook's append order.
$ws8 = $wb.Worksheets.Item("Personal_EL")
Set-Result $ws8 2 "Pass" "Fri Sep 29 17:28:43 EDT 2023"

# ---- Personal_IND (sheet6) -------------------------------------------
$ws6 = $wb.Worksheets.Item("Personal_IND")
Set-Result $ws6 2 "Pass" "Fri Sep 29 17:31:40 EDT 2023"
Set-Result $ws6 3 "Pass" "Fri Sep 29 17:32:23 EDT 2023"

# ---- Personal_JNT (sheet7) -------------------------------------------
$ws7 = $wb.Worksheets.Item("Personal_JNT")
Set-Result $ws7 2 "Pass" "Fri Sep 29 17:37:47 EDT 2023"
Set-Result $ws7 3 "Pass" "Fri Sep 29 17:38:37 EDT 2023"
Set-Result $ws7 4 "Fail" "Fri Sep 29 17:39:23 EDT 2023"
